$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at J, shifting Sector/Category/... etc one column to the right
$ws.Columns("J:J").Insert()

# Populate the new "Instrument" column
$ws.Range("J1").Value = "Instrument"
$ws.Range("J2").Value = "Stock"
$ws.Range("J3").Value = "Stock"

# Move the active selection to J4 (matches the post-edit cursor position)
$ws.Range("J4").Select()
